$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-03 Saturday" "2024-02-04 Sunday"

Replace-Text "320×9=" "276×6="
Replace-Text "660×5=" "447×3="
Replace-Text "833×9=" "721×4="
Replace-Text "817×5=" "533×7="
Replace-Text "351×2=" "837×6="
Replace-Text "937×9=" "155×2="
Replace-Text "217×6=" "354×8="
Replace-Text "667×9=" "942×6="
Replace-Text "427×3=" "343×9="
Replace-Text "983×5=" "960×2="
Replace-Text "298×7=" "923×7="
Replace-Text "795×4=" "915×9="
Replace-Text "986×3=" "540×5="
Replace-Text "484×2=" "273×5="
Replace-Text "217×2=" "209×6="
Replace-Text "186×2=" "301×4="
Replace-Text "495×6=" "115×4="
Replace-Text "770×4=" "701×8="
Replace-Text "559×2=" "607×5="
Replace-Text "696×5=" "494×5="
Replace-Text "621×7=" "769×4="
Replace-Text "350×9=" "898×6="
Replace-Text "994×9=" "921×2="
Replace-Text "377×4=" "661×9="
Replace-Text "455×8=" "363×5="
